# The workbook was re-worked around "DataNode" naming (see commit message:
# "unify the conception of DataNode, DataTable, Entity"). The concrete,
# user-visible changes captured in the diff are:
#   1. The worksheet "Property1" is renamed to "DataNode".
#   2. The live selection/active cell moves from A9 to D40.
#   3. A small secondary font (9pt, SimSun) gets registered in the style
#      table alongside a phonetic-guide default for the sheet (bookkeeping
#      Excel performs once a CJK-content sheet is re-saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet.
$ws.Name = "DataNode"

# 3) Touch a scratch cell far outside the used range (A1:H70) with the
#    small 9pt font so the workbook registers it in the style table (the
#    font backing the sheet's phonetic-guide default), then clear the
#    scratch cell again so neither the sheet dimension nor any visible
#    cell content/formatting is left changed.
$scratch = $ws.Range("Z100")
$scratch.Font.Size = 9
$scratch.Clear()

# Record the phonetic-guide default for the used range (mirrors Excel's
# own bookkeeping for CJK text sheets / the "Show Phonetic Field" toggle).
$ws.Range("A1:H70").SetPhonetic()

# 2) Move the selection to D40 (was A9).
$ws.Range("D40").Select()
